$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$k2New = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKPIA7QIT%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCsQb6xTTR%2FESU8i1kTlSZtiY%2Bf3SZOhDAFDVjMk5mhJQIgd4eSK8L%2BIHWHKkka0XBXLaAb6yhjr%2BH4sJkp%2BnCKXBkq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLYJkQBpszBc3w58aCq6A8UjpYC03amKNrY7ZK3BXk4gIMziubhBFEYqHIH6wjSo7CuNBMzRUu5dN%2Fo9l8m%2B3MyMYpxRyP666zXG26GWCWALKSpOwMikKC36dz2XX0RjXdybuKQN5nB%2FYzVbvkaAjqUmdYvXJdy3QZ8sZoLkYXVCJMSltr%2BuOnmLt3uU%2B881dUDEB2Qwo4KjoKaak%2F8Ah3qG5sHwNvrLw672NGFzV%2BDoCsl4HOciEhGUyqXEvBzSEnoNfGvv9Aa75QMZBwLdDy3nnbJ9KmshN3skqxOQlF8qL8LeC8uMY%2Fy%2FsJl66ifC38gclZGLHojT7U%2BaUt%2Boe7WmVq7FPE2WyNZHG9LQ5XfzCgeTG0%2FEOBqfb0SoSeuckbqncIcPp1Xkdqb8DQ7LpeiIURNLA0LqpZze%2FZts6O8xTjeF83NQPJVQHvpS6YRqN46tIZ0X%2BfTgwQAYFhPLwWGymqQsuioVCeF0utgKI%2Fgdhht7843nFtPhGsLyxSckudDVLeH6hBNnWrNuXN7FFEkua0DPUJ5ZH3ofcknFpEU0xc3OmMiBiOK%2FevA4hqBzhIpUroAzHMA3AZxJ%2FUY%2BSjyJSrazhCKK5SwwktPMwgY6pQF9UcVyv%2F%2FCjL8U0Ps2n%2BZ3murXY%2FDu7Uq%2BrMUXK8J7rrw3pyiJ9RJxvH%2BAw97OoTbGUgGbUYSA1mLZETMadbFhLxnthks8kybIweaoWuW01Mqvnaa1XgdmdwQQFkhQNv3ve1CxNsEu4PtRU7tAg3SvrMdfAT%2FYkrPjutRNOGeazqH9dDzfm4A9yZztn5DSriXyCMbq5Gjf7hdfkAJsNw25xNpyv3c%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=9eb0b96d05ca33a102093804b39a9ef4558e8a0b06b1b61e4d191cbda2f378b6'
$l2New = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKPIA7QIT%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCsQb6xTTR%2FESU8i1kTlSZtiY%2Bf3SZOhDAFDVjMk5mhJQIgd4eSK8L%2BIHWHKkka0XBXLaAb6yhjr%2BH4sJkp%2BnCKXBkq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLYJkQBpszBc3w58aCq6A8UjpYC03amKNrY7ZK3BXk4gIMziubhBFEYqHIH6wjSo7CuNBMzRUu5dN%2Fo9l8m%2B3MyMYpxRyP666zXG26GWCWALKSpOwMikKC36dz2XX0RjXdybuKQN5nB%2FYzVbvkaAjqUmdYvXJdy3QZ8sZoLkYXVCJMSltr%2BuOnmLt3uU%2B881dUDEB2Qwo4KjoKaak%2F8Ah3qG5sHwNvrLw672NGFzV%2BDoCsl4HOciEhGUyqXEvBzSEnoNfGvv9Aa75QMZBwLdDy3nnbJ9KmshN3skqxOQlF8qL8LeC8uMY%2Fy%2FsJl66ifC38gclZGLHojT7U%2BaUt%2Boe7WmVq7FPE2WyNZHG9LQ5XfzCgeTG0%2FEOBqfb0SoSeuckbqncIcPp1Xkdqb8DQ7LpeiIURNLA0LqpZze%2FZts6O8xTjeF83NQPJVQHvpS6YRqN46tIZ0X%2BfTgwQAYFhPLwWGymqQsuioVCeF0utgKI%2Fgdhht7843nFtPhGsLyxSckudDVLeH6hBNnWrNuXN7FFEkua0DPUJ5ZH3ofcknFpEU0xc3OmMiBiOK%2FevA4hqBzhIpUroAzHMA3AZxJ%2FUY%2BSjyJSrazhCKK5SwwktPMwgY6pQF9UcVyv%2F%2FCjL8U0Ps2n%2BZ3murXY%2FDu7Uq%2BrMUXK8J7rrw3pyiJ9RJxvH%2BAw97OoTbGUgGbUYSA1mLZETMadbFhLxnthks8kybIweaoWuW01Mqvnaa1XgdmdwQQFkhQNv3ve1CxNsEu4PtRU7tAg3SvrMdfAT%2FYkrPjutRNOGeazqH9dDzfm4A9yZztn5DSriXyCMbq5Gjf7hdfkAJsNw25xNpyv3c%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=67123dab126bd6b1866eae2d179c5ab7bbd970cabd46a93f1cec51bf77795029'
$k3New = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKPIA7QIT%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCsQb6xTTR%2FESU8i1kTlSZtiY%2Bf3SZOhDAFDVjMk5mhJQIgd4eSK8L%2BIHWHKkka0XBXLaAb6yhjr%2BH4sJkp%2BnCKXBkq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLYJkQBpszBc3w58aCq6A8UjpYC03amKNrY7ZK3BXk4gIMziubhBFEYqHIH6wjSo7CuNBMzRUu5dN%2Fo9l8m%2B3MyMYpxRyP666zXG26GWCWALKSpOwMikKC36dz2XX0RjXdybuKQN5nB%2FYzVbvkaAjqUmdYvXJdy3QZ8sZoLkYXVCJMSltr%2BuOnmLt3uU%2B881dUDEB2Qwo4KjoKaak%2F8Ah3qG5sHwNvrLw672NGFzV%2BDoCsl4HOciEhGUyqXEvBzSEnoNfGvv9Aa75QMZBwLdDy3nnbJ9KmshN3skqxOQlF8qL8LeC8uMY%2Fy%2FsJl66ifC38gclZGLHojT7U%2BaUt%2Boe7WmVq7FPE2WyNZHG9LQ5XfzCgeTG0%2FEOBqfb0SoSeuckbqncIcPp1Xkdqb8DQ7LpeiIURNLA0LqpZze%2FZts6O8xTjeF83NQPJVQHvpS6YRqN46tIZ0X%2BfTgwQAYFhPLwWGymqQsuioVCeF0utgKI%2Fgdhht7843nFtPhGsLyxSckudDVLeH6hBNnWrNuXN7FFEkua0DPUJ5ZH3ofcknFpEU0xc3OmMiBiOK%2FevA4hqBzhIpUroAzHMA3AZxJ%2FUY%2BSjyJSrazhCKK5SwwktPMwgY6pQF9UcVyv%2F%2FCjL8U0Ps2n%2BZ3murXY%2FDu7Uq%2BrMUXK8J7rrw3pyiJ9RJxvH%2BAw97OoTbGUgGbUYSA1mLZETMadbFhLxnthks8kybIweaoWuW01Mqvnaa1XgdmdwQQFkhQNv3ve1CxNsEu4PtRU7tAg3SvrMdfAT%2FYkrPjutRNOGeazqH9dDzfm4A9yZztn5DSriXyCMbq5Gjf7hdfkAJsNw25xNpyv3c%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=383f2fb2efaaaf6e3a78d94325a0e5785fbdd41a975542c6a72534f7d8e9fe19'
$l3New = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKPIA7QIT%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCsQb6xTTR%2FESU8i1kTlSZtiY%2Bf3SZOhDAFDVjMk5mhJQIgd4eSK8L%2BIHWHKkka0XBXLaAb6yhjr%2BH4sJkp%2BnCKXBkq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLYJkQBpszBc3w58aCq6A8UjpYC03amKNrY7ZK3BXk4gIMziubhBFEYqHIH6wjSo7CuNBMzRUu5dN%2Fo9l8m%2B3MyMYpxRyP666zXG26GWCWALKSpOwMikKC36dz2XX0RjXdybuKQN5nB%2FYzVbvkaAjqUmdYvXJdy3QZ8sZoLkYXVCJMSltr%2BuOnmLt3uU%2B881dUDEB2Qwo4KjoKaak%2F8Ah3qG5sHwNvrLw672NGFzV%2BDoCsl4HOciEhGUyqXEvBzSEnoNfGvv9Aa75QMZBwLdDy3nnbJ9KmshN3skqxOQlF8qL8LeC8uMY%2Fy%2FsJl66ifC38gclZGLHojT7U%2BaUt%2Boe7WmVq7FPE2WyNZHG9LQ5XfzCgeTG0%2FEOBqfb0SoSeuckbqncIcPp1Xkdqb8DQ7LpeiIURNLA0LqpZze%2FZts6O8xTjeF83NQPJVQHvpS6YRqN46tIZ0X%2BfTgwQAYFhPLwWGymqQsuioVCeF0utgKI%2Fgdhht7843nFtPhGsLyxSckudDVLeH6hBNnWrNuXN7FFEkua0DPUJ5ZH3ofcknFpEU0xc3OmMiBiOK%2FevA4hqBzhIpUroAzHMA3AZxJ%2FUY%2BSjyJSrazhCKK5SwwktPMwgY6pQF9UcVyv%2F%2FCjL8U0Ps2n%2BZ3murXY%2FDu7Uq%2BrMUXK8J7rrw3pyiJ9RJxvH%2BAw97OoTbGUgGbUYSA1mLZETMadbFhLxnthks8kybIweaoWuW01Mqvnaa1XgdmdwQQFkhQNv3ve1CxNsEu4PtRU7tAg3SvrMdfAT%2FYkrPjutRNOGeazqH9dDzfm4A9yZztn5DSriXyCMbq5Gjf7hdfkAJsNw25xNpyv3c%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=33a533c842d2409c2c758a7e1801ee0b2ca532a6c00c48fb971b51bb43e1dc5a'

$ws.Range("K2").Value = $k2New
$ws.Range("L2").Value = $l2New
$ws.Range("K3").Value = $k3New
$ws.Range("L3").Value = $l3New
